$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted at row 25 ("Fruta / hortaliza, semanal"
# commit): every existing data row from 25 downward shifts down by one,
# and a brand-new row of observations is written into the now-empty row 25.
$ws.Rows("25:25").Insert()

$ws.Range("A25").Value = 6
$ws.Range("B25").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C25").Value = "Metropolitana"
$ws.Range("D25").Value = 45125
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = 100112035
$ws.Range("G25").Value = "Bruselas (repollito)"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 580
$ws.Range("K25").Value = 17000
$ws.Range("L25").Value = 18000
$ws.Range("M25").Value = 17397
$ws.Range("N25").Value = "$/malla 15 kilos"
$ws.Range("O25").Value = "Provincia de Quillota"
$ws.Range("P25").Value = 1160
$ws.Range("Q25").Value = 15
$ws.Range("R25").Value = "Hortaliza"
